$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44176
$ws.Range("J2").Value = 1300
$ws.Range("K2").Value = 350
$ws.Range("L2").Value = 400
$ws.Range("M2").Value = 375
$ws.Range("P2").Value = 375

# Row 3
$ws.Range("D3").Value = 44231
$ws.Range("J3").Value = 200
$ws.Range("K3").Value = 180
$ws.Range("L3").Value = 200
$ws.Range("M3").Value = 190
$ws.Range("O3").Value = 'Región de Arica y Parinacota'
$ws.Range("P3").Value = 190

# Row 4
$ws.Range("D4").Value = 44175
$ws.Range("I4").Value = 'Segunda'
$ws.Range("K4").Value = 400
$ws.Range("L4").Value = 430
$ws.Range("M4").Value = 415
$ws.Range("P4").Value = 415

# Row 5
$ws.Range("D5").Value = 44224
$ws.Range("I5").Value = 'Segunda'
$ws.Range("J5").Value = 1200
$ws.Range("K5").Value = 230
$ws.Range("L5").Value = 250
$ws.Range("M5").Value = 240
$ws.Range("P5").Value = 240

# Row 6
$ws.Range("D6").Value = 44224
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 200
$ws.Range("L6").Value = 230
$ws.Range("M6").Value = 215
$ws.Range("O6").Value = 'Región de Arica y Parinacota'
$ws.Range("P6").Value = 215

# Row 7
$ws.Range("D7").Value = 44168
$ws.Range("J7").Value = 1700
$ws.Range("K7").Value = 430
$ws.Range("L7").Value = 450
$ws.Range("M7").Value = 440
$ws.Range("P7").Value = 440

# Row 8
$ws.Range("D8").Value = 44253
$ws.Range("I8").Value = 'Segunda'
$ws.Range("J8").Value = 1200
$ws.Range("K8").Value = 270
$ws.Range("L8").Value = 280
$ws.Range("M8").Value = 275
$ws.Range("P8").Value = 275

# Row 9
$ws.Range("D9").Value = 44214
$ws.Range("I9").Value = 'Primera'
$ws.Range("J9").Value = 1200
$ws.Range("K9").Value = 400
$ws.Range("L9").Value = 450
$ws.Range("M9").Value = 425
$ws.Range("P9").Value = 425

# Row 10
$ws.Range("D10").Value = 44547
$ws.Range("I10").Value = 'Primera'
$ws.Range("K10").Value = 350
$ws.Range("L10").Value = 370
$ws.Range("M10").Value = 360
$ws.Range("P10").Value = 360

# Row 11
$ws.Range("D11").Value = 44229
$ws.Range("K11").Value = 230
$ws.Range("L11").Value = 250
$ws.Range("M11").Value = 240
$ws.Range("O11").Value = 'Perú'
$ws.Range("P11").Value = 240

# Row 12
$ws.Range("D12").Value = 44659
$ws.Range("I12").Value = 'Primera'
$ws.Range("J12").Value = 1000
$ws.Range("K12").Value = 350
$ws.Range("L12").Value = 380
$ws.Range("M12").Value = 365
$ws.Range("P12").Value = 365

# Row 13
$ws.Range("D13").Value = 44172
$ws.Range("J13").Value = 1600
$ws.Range("K13").Value = 400
$ws.Range("L13").Value = 420
$ws.Range("M13").Value = 410
$ws.Range("O13").Value = 'Perú'
$ws.Range("P13").Value = 410

# Row 14
$ws.Range("D14").Value = 44580
$ws.Range("I14").Value = 'Primera'
$ws.Range("K14").Value = 380
$ws.Range("L14").Value = 400
$ws.Range("M14").Value = 390
$ws.Range("O14").Value = 'Región Metropolitana'
$ws.Range("P14").Value = 390

# Row 15
$ws.Range("D15").Value = 44530
$ws.Range("J15").Value = 1300
$ws.Range("K15").Value = 450
$ws.Range("L15").Value = 480
$ws.Range("M15").Value = 465
$ws.Range("O15").Value = 'Perú'
$ws.Range("P15").Value = 465

# Row 16
$ws.Range("D16").Value = 44162
$ws.Range("J16").Value = 900
$ws.Range("K16").Value = 500
$ws.Range("L16").Value = 550
$ws.Range("M16").Value = 525
$ws.Range("P16").Value = 525

# Row 17
$ws.Range("D17").Value = 44162
$ws.Range("J17").Value = 1200
$ws.Range("K17").Value = 500
$ws.Range("L17").Value = 550
$ws.Range("M17").Value = 525
$ws.Range("P17").Value = 525

# Row 18
$ws.Range("D18").Value = 44589
$ws.Range("I18").Value = 'Primera'
$ws.Range("J18").Value = 900
$ws.Range("K18").Value = 325
$ws.Range("M18").Value = 338
$ws.Range("O18").Value = 'Perú'
$ws.Range("P18").Value = 338

# Row 19
$ws.Range("D19").Value = 44202
$ws.Range("J19").Value = 1300
$ws.Range("K19").Value = 230
$ws.Range("L19").Value = 250
$ws.Range("M19").Value = 240
$ws.Range("P19").Value = 240

# Row 20
$ws.Range("D20").Value = 44201
$ws.Range("I20").Value = 'Segunda'
$ws.Range("J20").Value = 1800
$ws.Range("K20").Value = 250
$ws.Range("L20").Value = 270
$ws.Range("M20").Value = 260
$ws.Range("P20").Value = 260

# Row 21
$ws.Range("D21").Value = 44609
$ws.Range("I21").Value = 'Primera'
$ws.Range("J21").Value = 1200
$ws.Range("K21").Value = 280
$ws.Range("L21").Value = 300
$ws.Range("M21").Value = 290
$ws.Range("P21").Value = 290

# Row 22
$ws.Range("D22").Value = 44301
$ws.Range("I22").Value = 'Segunda'
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 280
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = 290
$ws.Range("P22").Value = 290

# Row 23
$ws.Range("D23").Value = 44217
$ws.Range("K23").Value = 300
$ws.Range("L23").Value = 350
$ws.Range("M23").Value = 325
$ws.Range("P23").Value = 325

# Row 24
$ws.Range("D24").Value = 44650
$ws.Range("I24").Value = 'Segunda'
$ws.Range("K24").Value = 325
$ws.Range("L24").Value = 350
$ws.Range("M24").Value = 338
$ws.Range("P24").Value = 338

# Row 25
$ws.Range("D25").Value = 44166
$ws.Range("J25").Value = 1700
$ws.Range("K25").Value = 500
$ws.Range("L25").Value = 530
$ws.Range("M25").Value = 515
$ws.Range("P25").Value = 515

# Row 26
$ws.Range("D26").Value = 44586
$ws.Range("I26").Value = 'Tercera'
$ws.Range("J26").Value = 500
$ws.Range("K26").Value = 330
$ws.Range("L26").Value = 350
$ws.Range("M26").Value = 340
$ws.Range("O26").Value = 'Región de Arica y Parinacota'
$ws.Range("P26").Value = 340

# Row 27
$ws.Range("D27").Value = 44160
$ws.Range("I27").Value = 'Segunda'
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 500
$ws.Range("L27").Value = 550
$ws.Range("M27").Value = 525
$ws.Range("P27").Value = 525

# Row 28
$ws.Range("D28").Value = 44243
$ws.Range("K28").Value = 300
$ws.Range("L28").Value = 320
$ws.Range("M28").Value = 310
$ws.Range("P28").Value = 310

# Row 29
$ws.Range("D29").Value = 44243
$ws.Range("J29").Value = 800
$ws.Range("K29").Value = 300
$ws.Range("L29").Value = 320
$ws.Range("M29").Value = 310
$ws.Range("P29").Value = 310

# Row 30
$ws.Range("D30").Value = 44251
$ws.Range("I30").Value = 'Primera'
$ws.Range("J30").Value = 1200
$ws.Range("K30").Value = 250
$ws.Range("L30").Value = 280
$ws.Range("M30").Value = 265
$ws.Range("P30").Value = 265

# Row 31
$ws.Range("D31").Value = 44603
$ws.Range("I31").Value = 'Tercera'
$ws.Range("J31").Value = 300
$ws.Range("K31").Value = 280
$ws.Range("L31").Value = 300
$ws.Range("M31").Value = 290
$ws.Range("O31").Value = 'Región de Arica y Parinacota'
$ws.Range("P31").Value = 290

# Row 32
$ws.Range("D32").Value = 44575
$ws.Range("J32").Value = 1200
$ws.Range("K32").Value = 380
$ws.Range("L32").Value = 400
$ws.Range("M32").Value = 390
$ws.Range("P32").Value = 390

# Row 33
$ws.Range("D33").Value = 44523
$ws.Range("J33").Value = 1000
$ws.Range("K33").Value = 550
$ws.Range("L33").Value = 580
$ws.Range("M33").Value = 565
$ws.Range("P33").Value = 565

# Row 34
$ws.Range("D34").Value = 44566
$ws.Range("J34").Value = 1000
$ws.Range("K34").Value = 300
$ws.Range("L34").Value = 320
$ws.Range("M34").Value = 310
$ws.Range("P34").Value = 310

# Row 35
$ws.Range("D35").Value = 44602
$ws.Range("I35").Value = 'Primera'
$ws.Range("J35").Value = 1300
$ws.Range("K35").Value = 350
$ws.Range("L35").Value = 380
$ws.Range("M35").Value = 365
$ws.Range("O35").Value = 'Perú'
$ws.Range("P35").Value = 365

# Row 36
$ws.Range("D36").Value = 44602
$ws.Range("J36").Value = 900
$ws.Range("L36").Value = 330
$ws.Range("M36").Value = 315
$ws.Range("P36").Value = 315
